# Enhanced layout of leader board and score board
# Adds 5 new ranking columns (a,B,c,d,e) and extends the score board
# downward with additional score blocks, mirroring the existing pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new columns F1:J1 -----------------------------------
$ws.Range("F1").Value2 = "a"
$ws.Range("G1").Value2 = "B"
$ws.Range("H1").Value2 = "c"
$ws.Range("I1").Value2 = "d"
$ws.Range("J1").Value2 = "e"
# Match header styling used by A1:E1
$ws.Range("A1").Copy() | Out-Null
$ws.Range("F1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Row 26: replace text scores with numeric scores ------------------
$arr = New-Object 'object[,]' 6,3
$arr[0,0] = -5;  $arr[0,1] = -6;  $arr[0,2] = 11
$arr[1,0] = -5;  $arr[1,1] = -6;  $arr[1,2] = 11
$arr[2,0] = -5;  $arr[2,1] = -3;  $arr[2,2] = 8
$arr[3,0] = -5;  $arr[3,1] = -3;  $arr[3,2] = 8
$arr[4,0] = -6;  $arr[4,1] = -3;  $arr[4,2] = 9
$arr[5,0] = 0;   $arr[5,1] = 0;   $arr[5,2] = 0
$ws.Range("C26:E31").Value2 = $arr

# --- Rows 32:44 (columns F:H) -----------------------------------------
$arr = New-Object 'object[,]' 13,3
$arr[0,0]  = -5; $arr[0,1]  = -6; $arr[0,2]  = 11
$arr[1,0]  = 0;  $arr[1,1]  = 0;  $arr[1,2]  = 0
$arr[2,0]  = -5; $arr[2,1]  = -6; $arr[2,2]  = 11
$arr[3,0]  = -5; $arr[3,1]  = -6; $arr[3,2]  = 11
$arr[4,0]  = -5; $arr[4,1]  = -3; $arr[4,2]  = -2
$arr[5,0]  = -5; $arr[5,1]  = -6; $arr[5,2]  = 11
$arr[6,0]  = -6; $arr[6,1]  = -3; $arr[6,2]  = 9
$arr[7,0]  = -5; $arr[7,1]  = -3; $arr[7,2]  = 8
$arr[8,0]  = -5; $arr[8,1]  = -6; $arr[8,2]  = -3
$arr[9,0]  = 0;  $arr[9,1]  = 0;  $arr[9,2]  = 0
$arr[10,0] = -6; $arr[10,1] = -3; $arr[10,2] = 9
$arr[11,0] = -6; $arr[11,1] = -5; $arr[11,2] = 11
$arr[12,0] = -5; $arr[12,1] = -3; $arr[12,2] = -3
$ws.Range("F32:H44").Value2 = $arr

# --- Rows 45:55 (columns F:J) ------------------------------------------
$arr = New-Object 'object[,]' 11,5
$arr[0,0]  = -5; $arr[0,1]  = -6; $arr[0,2]  = -3; $arr[0,3]  = -3; $arr[0,4]  = 17
$arr[1,0]  = 0;  $arr[1,1]  = 0;  $arr[1,2]  = 0;  $arr[1,3]  = 0;  $arr[1,4]  = 0
$arr[2,0]  = 0;  $arr[2,1]  = 0;  $arr[2,2]  = 0;  $arr[2,3]  = 0;  $arr[2,4]  = 0
$arr[3,0]  = 0;  $arr[3,1]  = 0;  $arr[3,2]  = 0;  $arr[3,3]  = 0;  $arr[3,4]  = 0
$arr[4,0]  = -6; $arr[4,1]  = -3; $arr[4,2]  = -3; $arr[4,3]  = -3; $arr[4,4]  = 15
$arr[5,0]  = 0;  $arr[5,1]  = 0;  $arr[5,2]  = 0;  $arr[5,3]  = 0;  $arr[5,4]  = 0
$arr[6,0]  = -5; $arr[6,1]  = -3; $arr[6,2]  = 8;  $arr[6,3]  = 8;  $arr[6,4]  = 8
$arr[7,0]  = 0;  $arr[7,1]  = 0;  $arr[7,2]  = 0;  $arr[7,3]  = 0;  $arr[7,4]  = 0
$arr[8,0]  = 0;  $arr[8,1]  = 0;  $arr[8,2]  = 0;  $arr[8,3]  = 0;  $arr[8,4]  = 0
$arr[9,0]  = 0;  $arr[9,1]  = 0;  $arr[9,2]  = 0;  $arr[9,3]  = 0;  $arr[9,4]  = 0
$arr[10,0] = -5; $arr[10,1] = -3; $arr[10,2] = -3; $arr[10,3] = -2; $arr[10,4] = 13
$ws.Range("F45:J55").Value2 = $arr

# --- Rows 56:57 (columns C:D) -------------------------------------------
$arr = New-Object 'object[,]' 2,2
$arr[0,0] = -5; $arr[0,1] = -6
$arr[1,0] = -6; $arr[1,1] = 6
$ws.Range("C56:D57").Value2 = $arr

# --- Rows 58:61, column A -------------------------------------------------
$arr = New-Object 'object[,]' 4,1
$arr[0,0] = -5
$arr[1,0] = 18
$arr[2,0] = -10
$arr[3,0] = 0
$ws.Range("A58:A61").Value2 = $arr

# --- Rows 58:61, columns C:E ----------------------------------------------
$arr = New-Object 'object[,]' 4,3
$arr[0,0] = -3;  $arr[0,1] = -2;  $arr[0,2] = 0
$arr[1,0] = -9;  $arr[1,1] = -6;  $arr[1,2] = -3
$arr[2,0] = -10; $arr[2,1] = -10; $arr[2,2] = 30
$arr[3,0] = 0;   $arr[3,1] = 0;   $arr[3,2] = 0
$ws.Range("C58:E61").Value2 = $arr

# --- Row 62: totals stored as text "-1" (A62, C62, D62, E62) -------------
foreach ($addr in @("A62", "C62", "D62", "E62")) {
    $cell = $ws.Range($addr)
    $cell.Value2 = "'-1"
    $cell.Style = "Normal"
}
